$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text so that
# values like "1.009" or "27.044.82" are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.044.82"
$ws.Range("E2").Value = "  -0.90%  "

# Row 3
$ws.Range("D3").Value = "1.830.77"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "311.52"
$ws.Range("E5").Value = "  -0.79%  "

# Row 6
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("D7").Value = "0.4638"
$ws.Range("E7").Value = "  -2.12%  "

# Row 8
$ws.Range("D8").Value = "0.3717"
$ws.Range("E8").Value = "  +1.01%  "

# Row 9
$ws.Range("D9").Value = "0.07377"
$ws.Range("E9").Value = "  -0.86%  "

# Row 10
$ws.Range("D10").Value = "0.8660"
$ws.Range("E10").Value = "  -2.09%  "

# Row 11
$ws.Range("D11").Value = "20.02"
$ws.Range("E11").Value = "  -2.23%  "

# Row 12
$ws.Range("D12").Value = "0.07838"
$ws.Range("E12").Value = "  +7.39%  "

# Row 13
$ws.Range("D13").Value = "1.836.90"
$ws.Range("E13").Value = "  -3.97%  "

# Row 14
$ws.Range("D14").Value = "6.628"
$ws.Range("E14").Value = "  +1.20%  "

# Row 15
$ws.Range("D15").Value = "5.368"
$ws.Range("E15").Value = "  -0.92%  "

# Row 16
$ws.Range("D16").Value = "92.09"
$ws.Range("E16").Value = "  -1.75%  "

# Row 17
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("D18").Value = "0.000008964"
$ws.Range("E18").Value = "  +2.07%  "

# Row 20
$ws.Range("E20").Value = "  -0.61%  "

# Row 21
$ws.Range("D21").Value = "27.081.98"
$ws.Range("E21").Value = "  -2.17%  "

# Row 22
$ws.Range("D22").Value = "5.172"
$ws.Range("E22").Value = "  -2.08%  "

# Row 23
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  -0.55%  "

# Row 24
$ws.Range("D24").Value = "2.066.59"
$ws.Range("E24").Value = "  -2.30%  "

# Row 25
$ws.Range("D25").Value = "152.96"
$ws.Range("E25").Value = "  +0.80%  "

# Row 26
$ws.Range("E26").Value = "  -2.85%  "

# Row 27
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  -2.11%  "

# Row 28
$ws.Range("D28").Value = "2.093"
$ws.Range("E28").Value = "  -1.64%  "

# Row 29
$ws.Range("D29").Value = "5.131"
$ws.Range("E29").Value = "  -1.79%  "

# Row 30
$ws.Range("D30").Value = "115.74"
$ws.Range("E30").Value = "  -1.16%  "

# Row 31
$ws.Range("D31").Value = "0.08871"
$ws.Range("E31").Value = "  -0.91%  "

# Row 32
$ws.Range("E32").Value = "  +0.84%  "

# Row 33
$ws.Range("D33").Value = "0.7310"
$ws.Range("E33").Value = "  -2.24%  "

# Row 34
$ws.Range("D34").Value = "4.451"
$ws.Range("E34").Value = "  -1.85%  "

# Row 35
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -3.21%  "

# Row 36
$ws.Range("D36").Value = "2.471"
$ws.Range("E36").Value = "  +2.08%  "

# Row 37
$ws.Range("D37").Value = "1.079"
$ws.Range("E37").Value = "  -1.49%  "

# Row 38
$ws.Range("D38").Value = "0.01950"
$ws.Range("E38").Value = "  -0.15%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "7.383"
$ws.Range("E39").Value = "  +2.18%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05242"
$ws.Range("E40").Value = "  -1.72%  "

# Row 41
$ws.Range("D41").Value = "2.934"
$ws.Range("E41").Value = "  -0.73%  "

# Row 42
$ws.Range("D42").Value = "0.5175"
$ws.Range("E42").Value = "  -2.10%  "

# Row 43
$ws.Range("D43").Value = "0.1633"
$ws.Range("E43").Value = "  -1.44%  "

# Row 44
$ws.Range("D44").Value = "0.8573"
$ws.Range("E44").Value = "  -15.07%  "

# Row 45
$ws.Range("D45").Value = "8.233"
$ws.Range("E45").Value = "  -2.97%  "

# Row 46
$ws.Range("D46").Value = "0.4838"
$ws.Range("E46").Value = "  -1.39%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.24"
$ws.Range("E47").Value = "  -2.28%  "

# Row 48
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("D49").Value = "102.84"
$ws.Range("E49").Value = "  -2.08%  "

# Row 50
$ws.Range("D50").Value = "1.625"
$ws.Range("E50").Value = "  -2.19%  "

# Row 51
$ws.Range("D51").Value = "0.06244"
$ws.Range("E51").Value = "  -0.87%  "
